$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: __init__ / Attributes are set to parameter values.
$ws.Range("D7").Value = "Attributes are set to parameter values."
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "account_number=2121 client_number=2222 balance=1000.00 date_created=today() minimum_balance=50.00"

# Row 8: __init__ / minimum_balance has invalid type.
$ws.Range("D8").Value = "minimum_balance has invalid type."
$ws.Range("F8").Value = 'account_number=2121 client_number=2222 balance=1000.00 date_created=today() minimum_balance="not_Mark"'

$ws.Range("G7").Value = "Object created"
$ws.Range("E8").Value = "None"
$ws.Range("G8").Value = 50

# Row 9: get_service_charges / balance greater than minimum balance
$ws.Range("D9").Value = "balance greater than minimum balance"
$ws.Range("E9").Value = "minimum_balance=50.00"
$ws.Range("F9").Value = "balance=500.00"
$ws.Range("G9").Value = "BASE_SERVICE_CHARGE=0.50"

# Row 10: get_service_charges / balance  equal to minimum balance
$ws.Range("D10").Value = "balance  equal to minimum balance"
$ws.Range("E10").Value = "minimum_balance=50.00"
$ws.Range("F10").Value = "balance=50.00"
$ws.Range("G10").Value = "BASE_SERVICE_CHARGE=0.50"

# Row 11: get_service_charges / balance less than minimum balance
$ws.Range("D11").Value = "balance less than minimum balance"
$ws.Range("E11").Value = "minimum_balance=50.00"
$ws.Range("F11").Value = "balance=49.99"
$ws.Range("G11").Value = "BSC * SCP = 1.00"

# Developer name (C3) - replaces placeholder "Student Name"
$ws.Range("C3").Value = "Ralph Vitug"

# Row 12: __str__ / appropriate value returned based on attribute values.
$ws.Range("D12").Value = "appropriate value returned based on attribute values."
$ws.Range("E12").Value = "minimum_balance=50.00 balance=1000.00"
$ws.Range("F12").Value = "account_number=2121 client_number=2222 balance=1000.00 date_created=today() minimum_balance=50.00"
$ws.Range("G12").Value = "Account number: 2121                                                    Balance:$1,000.00                                                                     Minimum balance: $50                                                            Account type: Savings"

# Leave the view focused on the last cell touched, matching the final selection
$ws.Range("G12").Select() | Out-Null
